$wb = $excel.ActiveWorkbook

# --- Sheet: ALC (index 1) ---
$ws = $wb.Worksheets.Item(1)
$ws.Range("H6").Value = 175.125
$ws.Range("I6").Value = 85.85714
$ws.Range("J6").Value = 800
$ws.Range("K6").Value = 257.57142
$ws.Range("L6").Value = 2400
$ws.Range("M6").Value = -145.57142
$ws.Range("N6").Value = -2624
$ws.Range("H9").Value = 132.5
$ws.Range("I9").Value = 167.5
$ws.Range("J9").Value = 80
$ws.Range("K9").Value = 167.5
$ws.Range("L9").Value = 80
$ws.Range("M9").Value = 1.5
$ws.Range("N9").Value = -418
$ws.Range("H12").Value = 62501330
$ws.Range("I12").Value = 1429.8182
$ws.Range("J12").Value = 200001100
$ws.Range("K12").Value = 1429.8182
$ws.Range("L12").Value = 200001100
$ws.Range("M12").Value = -1259.8182
$ws.Range("N12").Value = -200001440
$ws.Range("H21").Value = 12473.308
$ws.Range("I21").Value = 14019.125
$ws.Range("J21").Value = 10000
$ws.Range("K21").Value = 14019.125
$ws.Range("L21").Value = 10000
$ws.Range("M21").Value = -13551.125
$ws.Range("N21").Value = -10936
$ws.Range("H23").Value = 12473.308
$ws.Range("I23").Value = 14019.125
$ws.Range("J23").Value = 10000
$ws.Range("K23").Value = 14019.125
$ws.Range("L23").Value = 10000
$ws.Range("M23").Value = -13785.125
$ws.Range("N23").Value = -10468
$ws.Range("H29").Value = 80
$ws.Range("J29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("N29").ClearContents()
$ws.Range("H38").Value = 250
$ws.Range("I38").Value = 250
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 750
$ws.Range("L38").Value = 0
$ws.Range("M38").Value = -378
$ws.Range("N38").ClearContents()
$ws.Range("H40").Value = 1767.5135
$ws.Range("I40").Value = 1554.4546
$ws.Range("J40").Value = 2080
$ws.Range("K40").Value = 1554.4546
$ws.Range("L40").Value = 2080
$ws.Range("M40").Value = -1379.4546
$ws.Range("N40").Value = -2430
$ws.Range("H58").Value = 855.13336
$ws.Range("I58").Value = 855.13336
$ws.Range("K58").Value = 2565.40008
$ws.Range("M58").Value = -2415.40008
$ws.Range("H74").Value = 4966
$ws.Range("I74").Value = 4932.6665
$ws.Range("J74").Value = 4999.3335
$ws.Range("K74").Value = 4932.6665
$ws.Range("L74").Value = 4999.3335
$ws.Range("M74").Value = -3996.6665
$ws.Range("N74").Value = -6871.3335
$ws.Range("H77").Value = 4966
$ws.Range("I77").Value = 4932.6665
$ws.Range("J77").Value = 4999.3335
$ws.Range("K77").Value = 24663.3325
$ws.Range("L77").Value = 24996.6675
$ws.Range("M77").Value = -19983.3325
$ws.Range("N77").Value = -34356.6675
$ws.Range("H87").Value = 29499.75
$ws.Range("J87").Value = 29499.75
$ws.Range("L87").Value = 29499.75
$ws.Range("N87").Value = -31995.75
$ws.Range("H90").Value = 29499.75
$ws.Range("J90").Value = 29499.75
$ws.Range("L90").Value = 88499.25
$ws.Range("N90").Value = -100979.25
$ws.Range("H111").Value = 3145.7778
$ws.Range("I111").Value = 2129
$ws.Range("J111").Value = 3436.2856
$ws.Range("K111").Value = 6387
$ws.Range("L111").Value = 10308.8568
$ws.Range("M111").Value = -3320
$ws.Range("N111").Value = -16442.8568
$ws.Range("H131").Value = 43481412
$ws.Range("I131").Value = 111112170
$ws.Range("K131").Value = 333336510
$ws.Range("M131").Value = -333331470

# --- Sheet: ARM (index 2) ---
$ws = $wb.Worksheets.Item(2)
$ws.Range("H32").Value = 4765.8535
$ws.Range("I32").Value = 4561.9355
$ws.Range("J32").Value = 5398
$ws.Range("K32").Value = 4561.9355
$ws.Range("L32").Value = 5398
$ws.Range("M32").Value = -4274.9355
$ws.Range("N32").Value = -5972
$ws.Range("H74").Value = 62202.445
$ws.Range("I74").Value = 72583.57000000001
$ws.Range("J74").Value = 25868.5
$ws.Range("K74").Value = 72583.57000000001
$ws.Range("L74").Value = 25868.5
$ws.Range("M74").Value = -71709.57000000001
$ws.Range("N74").Value = -27616.5
$ws.Range("H77").Value = 62202.445
$ws.Range("I77").Value = 72583.57000000001
$ws.Range("J77").Value = 25868.5
$ws.Range("K77").Value = 362917.85
$ws.Range("L77").Value = 129342.5
$ws.Range("M77").Value = -358549.85
$ws.Range("N77").Value = -138078.5

# --- Sheet: CUL (index 5) ---
$ws = $wb.Worksheets.Item(5)
$ws.Range("H33").Value = 148.3
$ws.Range("I33").Value = 177.28572
$ws.Range("J33").Value = 80.666664
$ws.Range("K33").Value = 1063.71432
$ws.Range("L33").Value = 483.999984
$ws.Range("M33").Value = -780.71432
$ws.Range("N33").Value = -1049.999984
$ws.Range("H109").Value = 1970.7693
$ws.Range("I109").Value = 980.3077
$ws.Range("J109").Value = 2961.2307
$ws.Range("K109").Value = 2940.9231
$ws.Range("L109").Value = 8883.6921
$ws.Range("M109").Value = -1900.9231
$ws.Range("N109").Value = -10963.6921
$ws.Range("H113").Value = 671.3125
$ws.Range("I113").Value = 791.5
$ws.Range("J113").Value = 599.2
$ws.Range("K113").Value = 2374.5
$ws.Range("L113").Value = 1797.6
$ws.Range("M113").Value = -204.5
$ws.Range("N113").Value = -6137.6
$ws.Range("H122").Value = 1311.5883
$ws.Range("J122").Value = 1600.3334
$ws.Range("L122").Value = 14403.0006
$ws.Range("N122").Value = -19303.0006

# --- Sheet: GSM (index 6) ---
$ws = $wb.Worksheets.Item(6)
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()

# --- Sheet: LTW (index 7) ---
$ws = $wb.Worksheets.Item(7)
$ws.Range("H61").Value = 2695.7334
$ws.Range("I61").Value = 2217
$ws.Range("J61").Value = 3242.8572
$ws.Range("K61").Value = 2217
$ws.Range("L61").Value = 3242.8572
$ws.Range("M61").Value = -2015
$ws.Range("N61").Value = -3646.8572
$ws.Range("H82").Value = 3273
$ws.Range("I82").Value = 2000
$ws.Range("J82").Value = 3555.889
$ws.Range("K82").Value = 2000
$ws.Range("L82").Value = 3555.889
$ws.Range("M82").Value = -1639
$ws.Range("N82").Value = -4277.889
$ws.Range("H85").Value = 3273
$ws.Range("I85").Value = 2000
$ws.Range("J85").Value = 3555.889
$ws.Range("K85").Value = 2000
$ws.Range("L85").Value = 3555.889
$ws.Range("M85").Value = -752
$ws.Range("N85").Value = -6051.889
$ws.Range("H113").Value = 2695.7334
$ws.Range("I113").Value = 2217
$ws.Range("J113").Value = 3242.8572
$ws.Range("K113").Value = 2217
$ws.Range("L113").Value = 3242.8572
$ws.Range("M113").Value = -47
$ws.Range("N113").Value = -7582.8572
